$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new values look like plain numbers.
# Force them to remain text (matching the source data which stores
# prices as inline strings) by setting an explicit Text number format
# before writing the value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.72'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.38'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.30'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '312.86'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.08'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.79'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.401'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.74'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.11'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.78'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '263.83'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.580'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.81'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.42'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.98'

# D-column cells whose new values are naturally non-numeric text
# (thousands separated by extra dots), so no special handling is needed.
$ws.Range("D2").Value = '57.720.73'
$ws.Range("D3").Value = '2.442.79'
$ws.Range("D9").Value = '2.441.45'
$ws.Range("D14").Value = '2.876.62'
$ws.Range("D15").Value = '57.650.57'
$ws.Range("D18").Value = '2.442.43'

# E-column percentage change cells (plain text, padded with spaces).
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("E6").Value = '  -1.30%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("E12").Value = '  -4.02%  '
$ws.Range("E13").Value = '  -2.81%  '
$ws.Range("E14").Value = '  -1.92%  '
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("E18").Value = '  -2.87%  '
$ws.Range("E19").Value = '  -2.80%  '
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("E29").Value = '  +3.91%  '
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("E33").Value = '  -4.30%  '
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("E37").Value = '  -4.83%  '
$ws.Range("E38").Value = '  -4.82%  '
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("E43").Value = '  -5.14%  '
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("E51").Value = '  -2.57%  '
